$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

# Fill previously-blank cells in the existing log rows with the literal
# text "nan" (matches the source data's missing-value convention).
$ws.Range("D2:O2").Value = "nan"
$ws.Range("D3,G3:K3,M3:O3").Value = "nan"
$ws.Range("D4:O4").Value = "nan"
$ws.Range("D5,H5,J5:K5,M5:O5").Value = "nan"
$ws.Range("E6:G6,I6:K6,M6:O6").Value = "nan"
$ws.Range("E7,G7:J7,M7:O7").Value = "nan"
$ws.Range("F8:K8,M8:N8").Value = "nan"
$ws.Range("E9,H9:K9,M9:O9").Value = "nan"
$ws.Range("D10:O10").Value = "nan"
$ws.Range("D11:O11").Value = "nan"
$ws.Range("D12:O12").Value = "nan"
$ws.Range("D13:O13").Value = "nan"
$ws.Range("B14:K14").Value = "nan"
$ws.Range("B15:K15").Value = "nan"
$ws.Range("B16:K16").Value = "nan"

# Append the new service-event row (row 17) for Card8.
$ws.Range("A17").Value = "'8"
$ws.Range("L17").Value = "11\2\2025"
$ws.Range("M17").Value = "تم تركيب مساحه خروج شريط(35*26*1)"
$ws.Range("N17").Value = "كسر مساحه خروج شريط"
$ws.Range("O17").Value = "فني"
